$wb = $excel.ActiveWorkbook

# Sheet1: Device - Uptime value change
$wsDevice = $wb.Worksheets.Item("Device")
$wsDevice.Range("G2").Value = "23 hours, 29 minutes"

# Sheet2: Mem_CPU - value changes
$wsMem = $wb.Worksheets.Item("Mem_CPU")
$wsMem.Range("C3").Value = 405609916
$wsMem.Range("D3").Value = 1282750868
$wsMem.Range("I3").NumberFormat = "@"
$wsMem.Range("I3").Value = "1%"
$wsMem.Range("I3").Style = "Normal"

# Sheet3: Buffer - value changes
$wsBuffer = $wb.Worksheets.Item("Buffer")
$wsBuffer.Range("C2").Value = 357556
$wsBuffer.Range("C3").Value = 348792
$wsBuffer.Range("C4").Value = 78291
$wsBuffer.Range("C5").Value = 50819
